$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.104507
$ws.Range("H2").Value = 6.313521
$ws.Range("I2").Value = 0.2863228602749089
$ws.Range("J2").Value = 0.2863228602749089
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.593702666666666
$ws.Range("N2").Value = 16.781108
$ws.Range("O2").Value = 0.1866080739693142
$ws.Range("P2").Value = 0.1866080739693142
$ws.Range("Q2").Value = 11.77198641791867
$ws.Range("R2").Value = 105.947877761268
$ws.Range("S2").Value = 0.05343015748928583
$ws.Range("T2").Value = 0.05343015748928583

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.104507
$ws.Range("H3").Value = 6.313521
$ws.Range("I3").Value = 0.2863228602749089
$ws.Range("J3").Value = 0.2863228602749089
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.165799333333334
$ws.Range("N3").Value = 9.497398
$ws.Range("O3").Value = 0.105612284272291
$ws.Range("P3").Value = 0.105612284272291
$ws.Range("Q3").Value = 6.662446857595334
$ws.Range("R3").Value = 59.962021718358
$ws.Range("S3").Value = 0.03023921131300915
$ws.Range("T3").Value = 0.03023921131300914

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.104507
$ws.Range("H4").Value = 6.313521
$ws.Range("I4").Value = 0.2863228602749089
$ws.Range("J4").Value = 0.2863228602749089
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 20.201379
$ws.Range("N4").Value = 60.604137
$ws.Range("O4").Value = 0.673925778926067
$ws.Range("P4").Value = 0.673925778926067
$ws.Range("Q4").Value = 42.513943515153
$ws.Range("R4").Value = 382.625491636377
$ws.Range("S4").Value = 0.1929603566351074
$ws.Range("T4").Value = 0.1929603566351074

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.104507
$ws.Range("H5").Value = 6.313521
$ws.Range("I5").Value = 0.2863228602749089
$ws.Range("J5").Value = 0.2863228602749089
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.014792333333333
$ws.Range("N5").Value = 3.044377
$ws.Range("O5").Value = 0.03385386283232782
$ws.Range("P5").Value = 0.03385386283232782
$ws.Range("Q5").Value = 2.135637569046333
$ws.Range("R5").Value = 19.220738121417
$ws.Range("S5").Value = 0.00969313483750653
$ws.Range("T5").Value = 0.00969313483750653

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.935596333333333
$ws.Range("H6").Value = 5.806789
$ws.Range("I6").Value = 0.2633421882168251
$ws.Range("J6").Value = 0.2633421882168251
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.593702666666666
$ws.Range("N6").Value = 16.781108
$ws.Range("O6").Value = 0.1866080739693142
$ws.Range("P6").Value = 0.1866080739693142
$ws.Range("Q6").Value = 10.82715037135689
$ws.Range("R6").Value = 97.444353342212
$ws.Range("S6").Value = 0.04914177853800637
$ws.Range("T6").Value = 0.04914177853800637

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.935596333333333
$ws.Range("H7").Value = 5.806789
$ws.Range("I7").Value = 0.2633421882168251
$ws.Range("J7").Value = 0.2633421882168251
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.165799333333334
$ws.Range("N7").Value = 9.497398
$ws.Range("O7").Value = 0.105612284272291
$ws.Range("P7").Value = 0.105612284272291
$ws.Range("Q7").Value = 6.127709581669112
$ws.Range("R7").Value = 55.14938623502201
$ws.Range("S7").Value = 0.02781217004284251
$ws.Range("T7").Value = 0.0278121700428425

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.935596333333333
$ws.Range("H8").Value = 5.806789
$ws.Range("I8").Value = 0.2633421882168251
$ws.Range("J8").Value = 0.2633421882168251
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 20.201379
$ws.Range("N8").Value = 60.604137
$ws.Range("O8").Value = 0.673925778926067
$ws.Range("P8").Value = 0.673925778926067
$ws.Range("Q8").Value = 39.101715120677
$ws.Range("R8").Value = 351.915436086093
$ws.Range("S8").Value = 0.1774730893181188
$ws.Range("T8").Value = 0.1774730893181188

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.935596333333333
$ws.Range("H9").Value = 5.806789
$ws.Range("I9").Value = 0.2633421882168251
$ws.Range("J9").Value = 0.2633421882168251
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.014792333333333
$ws.Range("N9").Value = 3.044377
$ws.Range("O9").Value = 0.03385386283232782
$ws.Range("P9").Value = 0.03385386283232782
$ws.Range("Q9").Value = 1.964228319494778
$ws.Range("R9").Value = 17.678054875453
$ws.Range("S9").Value = 0.008915150317857452
$ws.Range("T9").Value = 0.008915150317857452

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.414843
$ws.Range("H10").Value = 7.244529
$ws.Range("I10").Value = 0.3285447636310271
$ws.Range("J10").Value = 0.3285447636310271
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.593702666666666
$ws.Range("N10").Value = 16.781108
$ws.Range("O10").Value = 0.1866080739693142
$ws.Range("P10").Value = 0.1866080739693142
$ws.Range("Q10").Value = 13.50791372868133
$ws.Range("R10").Value = 121.571223558132
$ws.Range("S10").Value = 0.06130910555388956
$ws.Range("T10").Value = 0.06130910555388956

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.414843
$ws.Range("H11").Value = 7.244529
$ws.Range("I11").Value = 0.3285447636310271
$ws.Range("J11").Value = 0.3285447636310271
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 3.165799333333334
$ws.Range("N11").Value = 9.497398
$ws.Range("O11").Value = 0.105612284272291
$ws.Range("P11").Value = 0.105612284272291
$ws.Range("Q11").Value = 7.644908359504667
$ws.Range("R11").Value = 68.804175235542
$ws.Range("S11").Value = 0.0346983629727727
$ws.Range("T11").Value = 0.03469836297277269

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.414843
$ws.Range("H12").Value = 7.244529
$ws.Range("I12").Value = 0.3285447636310271
$ws.Range("J12").Value = 0.3285447636310271
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 20.201379
$ws.Range("N12").Value = 60.604137
$ws.Range("O12").Value = 0.673925778926067
$ws.Range("P12").Value = 0.673925778926067
$ws.Range("Q12").Value = 48.78315866849699
$ws.Range("R12").Value = 439.048428016473
$ws.Range("S12").Value = 0.2214147857421205
$ws.Range("T12").Value = 0.2214147857421205

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.414843
$ws.Range("H13").Value = 7.244529
$ws.Range("I13").Value = 0.3285447636310271
$ws.Range("J13").Value = 0.3285447636310271
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.014792333333333
$ws.Range("N13").Value = 3.044377
$ws.Range("O13").Value = 0.03385386283232782
$ws.Range("P13").Value = 0.03385386283232782
$ws.Range("Q13").Value = 2.450564162603666
$ws.Range("R13").Value = 22.055077463433
$ws.Range("S13").Value = 0.01112250936224436
$ws.Range("T13").Value = 0.01112250936224436

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.8951723333333333
$ws.Range("H14").Value = 2.685517
$ws.Range("I14").Value = 0.1217901878772388
$ws.Range("J14").Value = 0.1217901878772388
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 5.593702666666666
$ws.Range("N14").Value = 16.781108
$ws.Range("O14").Value = 0.1866080739693142
$ws.Range("P14").Value = 0.1866080739693142
$ws.Range("Q14").Value = 5.007327868092888
$ws.Range("R14").Value = 45.065950812836
$ws.Range("S14").Value = 0.02272703238813245
$ws.Range("T14").Value = 0.02272703238813245

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.8951723333333333
$ws.Range("H15").Value = 2.685517
$ws.Range("I15").Value = 0.1217901878772388
$ws.Range("J15").Value = 0.1217901878772388
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 3.165799333333334
$ws.Range("N15").Value = 9.497398
$ws.Range("O15").Value = 0.105612284272291
$ws.Range("P15").Value = 0.105612284272291
$ws.Range("Q15").Value = 2.833935976085111
$ws.Range("R15").Value = 25.505423784766
$ws.Range("S15").Value = 0.01286253994366668
$ws.Range("T15").Value = 0.01286253994366668

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.8951723333333333
$ws.Range("H16").Value = 2.685517
$ws.Range("I16").Value = 0.1217901878772388
$ws.Range("J16").Value = 0.1217901878772388
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 20.201379
$ws.Range("N16").Value = 60.604137
$ws.Range("O16").Value = 0.673925778926067
$ws.Range("P16").Value = 0.673925778926067
$ws.Range("Q16").Value = 18.083715575981
$ws.Range("R16").Value = 162.753440183829
$ws.Range("S16").Value = 0.0820775472307202
$ws.Range("T16").Value = 0.0820775472307202

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.8951723333333333
$ws.Range("H17").Value = 2.685517
$ws.Range("I17").Value = 0.1217901878772388
$ws.Range("J17").Value = 0.1217901878772388
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 1.014792333333333
$ws.Range("N17").Value = 3.044377
$ws.Range("O17").Value = 0.03385386283232782
$ws.Range("P17").Value = 0.03385386283232782
$ws.Range("Q17").Value = 0.9084140208787778
$ws.Range("R17").Value = 8.175726187909
$ws.Range("S17").Value = 0.004123068314719477
$ws.Range("T17").Value = 0.004123068314719477
